$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its literal text representation
# (values like "1.002" or "0.07240" must not become numbers)
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '28.891.34'
$ws.Range('E2').Value = '  +2.65%  '
$ws.Range('D3').Value = '1.879.87'
$ws.Range('E3').Value = '  +2.67%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '327.12'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = '0.4657'
$ws.Range('E7').Value = '  +0.99%  '
$ws.Range('D8').Value = '0.3947'
$ws.Range('E8').Value = '  +2.75%  '
$ws.Range('D9').Value = '0.07922'
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('D10').Value = '0.9770'
$ws.Range('E10').Value = '  +2.34%  '
$ws.Range('D11').Value = '22.41'
$ws.Range('E11').Value = '  +2.77%  '
$ws.Range('D12').Value = '1.900.74'
$ws.Range('E12').Value = '  +3.89%  '
$ws.Range('D13').Value = '5.764'
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('D14').Value = '6.968'
$ws.Range('E14').Value = '  +1.57%  '
$ws.Range('D15').Value = '0.07001'
$ws.Range('E15').Value = '  +2.23%  '
$ws.Range('D16').Value = '88.73'
$ws.Range('E16').Value = '  +2.69%  '
$ws.Range('D17').Value = '1.003'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '0.00001016'
$ws.Range('E18').Value = '  +2.62%  '
$ws.Range('D19').Value = '17.01'
$ws.Range('E19').Value = '  +1.14%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '28.869.79'
$ws.Range('E21').Value = '  +2.45%  '
$ws.Range('D22').Value = '5.360'
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('D23').Value = '11.12'
$ws.Range('E23').Value = '  +1.83%  '
$ws.Range('D24').Value = '2.113'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').Value = '2.066.65'
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('D27').Value = '19.42'
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('D28').Value = '5.790'
$ws.Range('E28').Value = '  +1.80%  '
$ws.Range('D29').Value = '2.013'
$ws.Range('E29').Value = '  +2.32%  '
$ws.Range('D30').Value = '119.93'
$ws.Range('E30').Value = '  +2.85%  '
$ws.Range('D31').Value = '0.09387'
$ws.Range('E31').Value = '  +1.41%  '
$ws.Range('D32').Value = '0.9451'
$ws.Range('E32').Value = '  +1.10%  '
$ws.Range('D33').Value = '5.329'
$ws.Range('E33').Value = '  +1.55%  '
$ws.Range('D34').Value = '1.355'
$ws.Range('E34').Value = '  +4.02%  '
$ws.Range('D35').Value = '3.340'
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('D36').Value = '0.05893'
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('D37').Value = '0.02120'
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').Value = '7.955'
$ws.Range('E38').Value = '  +4.79%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '1.147'
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('D40').Value = '0.5696'
$ws.Range('E40').Value = '  +1.99%  '
$ws.Range('D41').Value = '0.1792'
$ws.Range('E41').Value = '  +1.33%  '
$ws.Range('D42').Value = '9.979'
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('D43').Value = '0.07240'
$ws.Range('E43').Value = '  +3.31%  '
$ws.Range('D44').Value = '11.82'
$ws.Range('E44').Value = '  +2.35%  '
$ws.Range('D45').Value = '0.5339'
$ws.Range('E45').Value = '  +1.75%  '
$ws.Range('D46').Value = '1.135'
$ws.Range('E46').Value = '  -7.14%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = '2.120'
$ws.Range('E47').Value = '  -5.02%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '1.856'
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('D49').Value = '114.11'
$ws.Range('E49').Value = '  +1.95%  '
$ws.Range('D50').Value = '2.361'
$ws.Range('E50').Value = '  +1.56%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = '1.031'
$ws.Range('E51').Value = '  +2.48%  '
